# This script appends new sensor-log rows to several worksheets, matching
# the data produced by an automated logging run on 2026-01-30.
$wb = $excel.ActiveWorkbook

function Set-LogRow {
    param($ws, $rowNum, $date, $timestamp, $hour, $location, $value, $status)

    # Column A holds a literal date-like string ("2026-01-30"). Assigning it
    # directly would make Excel auto-convert it into a real date serial, so we
    # force text by leading with an apostrophe and then clear the resulting
    # quote-prefix style so the cell keeps the workbook default style.
    $cell = $ws.Cells.Item($rowNum, 1)
    $cell.Value = "'$date"
    $cell.Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = $timestamp
    $ws.Cells.Item($rowNum, 3).Value = $hour
    $ws.Cells.Item($rowNum, 4).Value = $location

    # Column E sometimes holds a percentage-looking string (e.g. "87.4%").
    # Same text-forcing trick is needed so it is not converted to a number.
    if ($value -match "%$") {
        $cell = $ws.Cells.Item($rowNum, 5)
        $cell.Value = "'$value"
        $cell.Style = "Normal"
    } else {
        $ws.Cells.Item($rowNum, 5).Value = $value
    }

    $ws.Cells.Item($rowNum, 6).Value = $status
}

# --- PIR sheet: append rows 184-207 ---
$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    ,@(184, '2026-01-30', '17:11:19', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(185, '2026-01-30', '17:11:19', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(186, '2026-01-30', '17:11:23', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(187, '2026-01-30', '17:11:28', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(188, '2026-01-30', '17:11:33', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(189, '2026-01-30', '17:11:38', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(190, '2026-01-30', '17:11:43', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(191, '2026-01-30', '17:11:48', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(192, '2026-01-30', '17:11:53', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(193, '2026-01-30', '17:11:58', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(194, '2026-01-30', '17:12:03', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(195, '2026-01-30', '17:12:08', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(196, '2026-01-30', '17:23:15', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(197, '2026-01-30', '17:23:17', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(198, '2026-01-30', '17:23:22', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(199, '2026-01-30', '17:23:27', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(200, '2026-01-30', '17:23:32', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(201, '2026-01-30', '17:23:37', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(202, '2026-01-30', '17:23:42', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(203, '2026-01-30', '17:23:47', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(204, '2026-01-30', '17:23:52', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(205, '2026-01-30', '17:23:57', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(206, '2026-01-30', '17:24:02', '17:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@(207, '2026-01-30', '17:24:07', '17:00', 'Bathroom', 'No Motion', 'Inactive')
)
foreach ($r in $rows) {
    Set-LogRow $ws $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}

# --- Humidity sheet: append rows 125-137 ---
$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    ,@(125, '2026-01-30', '17:11:19', '17:00', 'Bathroom', '87.4%', 'Active')
    ,@(126, '2026-01-30', '17:11:23', '17:00', 'Bathroom', '87.4%', 'Active')
    ,@(127, '2026-01-30', '17:11:28', '17:00', 'Bathroom', '86.4%', 'Active')
    ,@(128, '2026-01-30', '17:11:43', '17:00', 'Bathroom', '87.4%', 'Active')
    ,@(129, '2026-01-30', '17:11:53', '17:00', 'Bathroom', '87.3%', 'Active')
    ,@(130, '2026-01-30', '17:12:03', '17:00', 'Bathroom', '86.0%', 'Active')
    ,@(131, '2026-01-30', '17:23:16', '17:00', 'Bathroom', '87.5%', 'Active')
    ,@(132, '2026-01-30', '17:23:17', '17:00', 'Bathroom', '86.5%', 'Active')
    ,@(133, '2026-01-30', '17:23:27', '17:00', 'Bathroom', '87.4%', 'Active')
    ,@(134, '2026-01-30', '17:23:32', '17:00', 'Bathroom', '87.5%', 'Active')
    ,@(135, '2026-01-30', '17:23:37', '17:00', 'Bathroom', '87.5%', 'Active')
    ,@(136, '2026-01-30', '17:23:48', '17:00', 'Bathroom', '87.5%', 'Active')
    ,@(137, '2026-01-30', '17:23:53', '17:00', 'Bathroom', '87.5%', 'Active')
)
foreach ($r in $rows) {
    Set-LogRow $ws $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}

# --- Proximity sheet: append rows 49-51 ---
$ws = $wb.Worksheets.Item("Proximity")
$rows = @(
    ,@(49, '2026-01-30', '17:23:41', '17:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
    ,@(50, '2026-01-30', '17:23:43', '17:00', 'Living Room Main Door', 'EXIT', 'User EXITED Living Room Main Door')
    ,@(51, '2026-01-30', '17:23:49', '17:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
)
foreach ($r in $rows) {
    Set-LogRow $ws $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}

# --- mmWave sheet: append rows 42-47 ---
$ws = $wb.Worksheets.Item("mmWave")
$rows = @(
    ,@(42, '2026-01-30', '17:23:11', '17:00', 'Living Room', 'FALL_DETECTED', 'EMERGENCY')
    ,@(43, '2026-01-30', '17:23:14', '17:00', 'Living Room', 'FALL_DETECTED', 'EMERGENCY')
    ,@(44, '2026-01-30', '17:23:15', '17:00', 'Living Room', 'PRESENCE_DETECTED', 'Active')
    ,@(45, '2026-01-30', '17:23:22', '17:00', 'Living Room', 'PRESENCE_DETECTED', 'Active')
    ,@(46, '2026-01-30', '17:23:33', '17:00', 'Living Room', 'PRESENCE_DETECTED', 'Active')
    ,@(47, '2026-01-30', '17:23:47', '17:00', 'Living Room', 'PRESENCE_DETECTED', 'Active')
)
foreach ($r in $rows) {
    Set-LogRow $ws $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}

# --- Camera sheet: append rows 14-16 ---
$ws = $wb.Worksheets.Item("Camera")
$rows = @(
    ,@(14, '2026-01-30', '17:23:41', '17:00', 'Living Room Main Door', 'Image Captured (ENTER)', 'Active')
    ,@(15, '2026-01-30', '17:23:43', '17:00', 'Living Room Main Door', 'Image Captured (EXIT)', 'Active')
    ,@(16, '2026-01-30', '17:23:48', '17:00', 'Living Room Main Door', 'Image Captured (ENTER)', 'Active')
)
foreach ($r in $rows) {
    Set-LogRow $ws $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}
